$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "2505" period row for CARLOS ANDRES BALDIRIS REINO
# (row 16). Deleting shifts every following row up by one, which also
# takes care of moving the signature block from rows 29/30 to 28/29.
$ws.Rows.Item(16).Delete()

# Rebuild the worker/period detail table (now rows 16-23) in the new
# order: ERIKA's periods first (ascending 2206-2208), then CARLOS's
# remaining periods (ascending 2403-2406), then CAROLINA unchanged.

# ERIKA PATRICIA MATURANA ROSENSTAND
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047372612"
$ws.Range("D16").Value = "ERIKA PATRICIA MATURANA ROSENSTAND"
$ws.Range("E16").Value = "2206"
$ws.Range("F16").Value = 68000
$ws.Range("G16").Value = 1700000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047372612"
$ws.Range("D17").Value = "ERIKA PATRICIA MATURANA ROSENSTAND"
$ws.Range("E17").Value = "2207"
$ws.Range("F17").Value = 68000
$ws.Range("G17").Value = 1700000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047372612"
$ws.Range("D18").Value = "ERIKA PATRICIA MATURANA ROSENSTAND"
$ws.Range("E18").Value = "2208"
$ws.Range("F18").Value = 68000
$ws.Range("G18").Value = 1700000

# CARLOS ANDRES BALDIRIS REINO
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1047393841"
$ws.Range("D19").Value = "CARLOS ANDRES BALDIRIS REINO"
$ws.Range("E19").Value = "2403"
$ws.Range("F19").Value = 34666
$ws.Range("G19").Value = 1790932

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1047393841"
$ws.Range("D20").Value = "CARLOS ANDRES BALDIRIS REINO"
$ws.Range("E20").Value = "2404"
$ws.Range("F20").Value = 52000
$ws.Range("G20").Value = 1790932

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1047393841"
$ws.Range("D21").Value = "CARLOS ANDRES BALDIRIS REINO"
$ws.Range("E21").Value = "2405"
$ws.Range("F21").Value = 52000
$ws.Range("G21").Value = 1790932

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1047393841"
$ws.Range("D22").Value = "CARLOS ANDRES BALDIRIS REINO"
$ws.Range("E22").Value = "2406"
$ws.Range("F22").Value = 52000
$ws.Range("G22").Value = 1790932

# CAROLINA MARIMON SIMARRA (content unchanged, just shifted up a row)
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1001833348"
$ws.Range("D23").Value = "CAROLINA MARIMON SIMARRA"
$ws.Range("E23").Value = "2507"
$ws.Range("F23").Value = 68000
$ws.Range("G23").Value = 1700000

# Update summary figures: total overdue value and period count
$ws.Range("E11").Value = 462666
$ws.Range("F13").Value = 8
